# Apply updated GA / preprocessing results for EQ model rows 2-6
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 65.945913649036683
$ws.Range("D2").Value = 0.86533237141068364
$ws.Range("E2").Value = 66.843945765149201
$ws.Range("F2").Value = 0.97375581202086281
$ws.Range("G2").Value = 15.731999831838332

$ws.Range("C3").Value = 60.281550895109795
$ws.Range("D3").Value = 0.93453264185714435
$ws.Range("E3").Value = 49.272854016130395
$ws.Range("F3").Value = 1.2044502233804624
$ws.Range("G3").Value = 20.664665002975042

$ws.Range("B4").Value = 15
$ws.Range("C4").Value = 24.741887516206372
$ws.Range("D4").Value = 1.2863969302859308
$ws.Range("E4").Value = 35.177898398257199
$ws.Range("F4").Value = 1.3615393447435347
$ws.Range("G4").Value = 21.631647537449201

$ws.Range("B5").Value = 1939
$ws.Range("C5").Value = 75.879453113362416
$ws.Range("D5").Value = 0.72826914334596071
$ws.Range("E5").Value = 54.544833959881636
$ws.Range("F5").Value = 1.1675916874235441
$ws.Range("G5").Value = 20.385222085635686

$ws.Range("B6").Value = 38
$ws.Range("C6").Value = 82.034676572744075
$ws.Range("D6").Value = 0.62851530004341649
$ws.Range("E6").Value = 60.509269546238301
$ws.Range("F6").Value = 1.0627135685394753
$ws.Range("G6").Value = 17.355022868323601
